$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Write header + data rows (Set Pair / Word Pair) ---
$ws.Cells.Item(1, 1).Value = 'Set Pair'
$ws.Cells.Item(1, 2).Value = 'Word Pair'
$ws.Cells.Item(2, 1).Value = 'Bajta - Britto_2017'
$ws.Cells.Item(2, 2).Value = 'risk - risk level - hardware - technical factors - platform support level - design - architecture - maintainability - portability level - installability level - maintainability level - reliability level - reusability level - size report - team size - performance - time efficiency level - process efficiency level - quality level - availability - availability level - individual - personality - security - security level - reliability - robustness level - testing - testability level - effort hours - maintenance - no of team members - team capability - work team level'
$ws.Cells.Item(3, 1).Value = 'Bajta - Britto_2016'
$ws.Cells.Item(3, 2).Value = 'estimate value - estimator & provider - estimator - temporal distance - relationship.geographic distance - site.temporal distance - relationship.temporal distance - geographical distance - site.geographic distance'
$ws.Cells.Item(4, 1).Value = 'Bajta - Dasthi'
$ws.Cells.Item(4, 2).Value = 'constructive cost model - fuzzy similar - fuzzy logic - expert judgment - machine learning - artificial neural networks'
$ws.Cells.Item(5, 1).Value = 'Bajta - Mendes'
$ws.Cells.Item(5, 2).Value = 'not considered - validated theoretically - size report - early size metric - late size metric - implementation - functionality - performance - considered'
$ws.Cells.Item(6, 1).Value = 'Bajta - Usman'
$ws.Cells.Item(6, 2).Value = 'estimate value - estimation entity.other - number of entities estimated.value - estimation techniques.other - estimate value(s) - actual effort.value - not considered - not used - considered - not applicable - agile - customized scrum - scrum - group-based estimation - telecommunication - communications industry - near offshore - distributed: far offshore - distributed: distant onshore - distributed: near offshore - distributed: close onshore - design - maintainability - availability - reliability - maintenance - healthcare - health - size report - size.other - implementation - performance - analysis - execution - task - finance - financial - statistics analysis - individual - single - security - other - project domain.other - unit.other - se - value - system investigation - far offshore - expert judgment - expert judgement - close onshore - testing - effort hours - ideal hours - effort estimate.type.other - hours/days - distant onshore - no of team members - no. of team members'
$ws.Cells.Item(7, 1).Value = 'Britto_2017 - Bajta'
$ws.Cells.Item(7, 2).Value = 'technical factors - hardware - portability level - maintainability - risk level - risk - team capability - no of team members - time efficiency level - performance - effort hours - work team level - architecture - design - process efficiency level - availability level - availability - testability level - testing - personality - individual - robustness level - reliability - installability level - team size - size report - security level - security - quality level - maintainability level - maintenance - platform support level - reliability level - reusability level'
$ws.Cells.Item(8, 1).Value = 'Britto_2017 - Dasthi'
$ws.Cells.Item(8, 2).Value = 'software development experience - software life cycle management'
$ws.Cells.Item(9, 1).Value = 'Britto_2017 - Mendes'
$ws.Cells.Item(9, 2).Value = 'interface complexity - complexity - control flow complexity - media count - media - web page allocation - web application - project.type - model dependency.specific - adaptation complexity - program count - program/script - class complexity - class.length - cyclomatic complexity - web objects - web software application - web hypermedia application - model collection complexity - difficulty level - data flow complexity - motivation level - motivation - data usage complexity - new media count - media allocation - cohesion complexity - page complexity - output complexity - new complexity - input complexity - component complexity - model association complexity - total complexity - media duration - layout complexity'
$ws.Cells.Item(10, 1).Value = 'Britto_2017 - Usman'
$ws.Cells.Item(10, 2).Value = 'entity count - estimation entity.other - number of entities estimated.value - portability level - maintainability - project.infrastructure - project domain.other - team capability - no. of team members - project.type - effort estimate.type.other - time efficiency level - performance - work team level - time restriction - hours/days - use case count - use case - use case points method - not used - user case points - requirements novelty level - non functional requirements.other - data web points - point - international function point users group - function points - architecture - design - requirements clarity level - process efficiency level - availability level - availability - testability level - testing - processing requirements - object-oriented function points - requirements volatility level - robustness level - reliability - installability level - implementation - team size - security level - security - quality level - accuracy level.value - object-oriented heuristic function points - maintainability level - maintenance - reliability level - reusability level'
$ws.Cells.Item(11, 1).Value = 'Britto_2016 - Bajta'
$ws.Cells.Item(11, 2).Value = 'estimator & provider - estimate value - relationship.geographic distance - temporal distance - geographical distance - site.geographic distance - estimator - site.temporal distance - relationship.temporal distance'
$ws.Cells.Item(12, 1).Value = 'Britto_2016 - Usman'
$ws.Cells.Item(12, 2).Value = 'semi-distributed - distributed: far offshore - distributed: near offshore - distribution - estimator & provider - estimation entity.other - estimation techniques.other - estimate value(s) - relationship.location - co-located - centralized - estimator - distributed - distributed: distant onshore - distributed: close onshore'
$ws.Cells.Item(13, 1).Value = 'Dasthi - Bajta'
$ws.Cells.Item(13, 2).Value = 'expert judgment - constructive cost model - artificial neural networks - machine learning - fuzzy logic - fuzzy similar'
$ws.Cells.Item(14, 1).Value = 'Dasthi - Britto_2017'
$ws.Cells.Item(14, 2).Value = 'software life cycle management - software development experience'
$ws.Cells.Item(15, 1).Value = 'Dasthi - Usman'
$ws.Cells.Item(15, 2).Value = 'expert judgment - expert judgement - analogy-based - analogy'
$ws.Cells.Item(16, 1).Value = 'Mendes - Bajta'
$ws.Cells.Item(16, 2).Value = 'validated theoretically - not considered - considered - functionality - implementation - performance - early size metric - size report - late size metric'
$ws.Cells.Item(17, 1).Value = 'Mendes - Britto_2017'
$ws.Cells.Item(17, 2).Value = 'web software application - web objects - motivation - motivation level - program/script - program count - web hypermedia application - complexity - interface complexity - control flow complexity - adaptation complexity - class complexity - cyclomatic complexity - model collection complexity - difficulty level - data flow complexity - data usage complexity - cohesion complexity - page complexity - output complexity - new complexity - input complexity - component complexity - model association complexity - total complexity - layout complexity - model dependency.specific - project.type - class.length - web application - web page allocation - media - media count - new media count - media allocation - media duration'
$ws.Cells.Item(18, 1).Value = 'Mendes - Usman'
$ws.Cells.Item(18, 2).Value = 'validated theoretically - not considered - considered - functionality - use case - maintainability - task - implementation - non functional requirements.other - performance - solution-oriented metric - considered without any metric - early size metric - size.other - problem-oriented metric - late size metric'
$ws.Cells.Item(19, 1).Value = 'Usman - Bajta'
$ws.Cells.Item(19, 2).Value = 'ideal hours - effort hours - size.other - size report - other - effort estimate.type.other - no. of team members - no of team members - not considered - considered - se - distributed: far offshore - near offshore - far offshore - distant onshore - distributed: distant onshore - close onshore - customized scrum - agile - maintainability - availability - reliability - maintenance - design - hours/days - task - execution - distributed: near offshore - communications industry - telecommunication - implementation - performance - project domain.other - estimation entity.other - estimate value - analysis - number of entities estimated.value - estimation techniques.other - group-based estimation - not used - statistics analysis - system investigation - distributed: close onshore - security - unit.other - expert judgement - expert judgment - estimate value(s) - value - single - individual - health - healthcare - financial - finance - not applicable - testing - scrum - actual effort.value'
$ws.Cells.Item(20, 1).Value = 'Usman - Britto_2017'
$ws.Cells.Item(20, 2).Value = 'effort estimate.type.other - project.type - use case - use case count - no. of team members - team capability - work team level - team size - use case points method - function points - international function point users group - object-oriented function points - object-oriented heuristic function points - point - data web points - maintainability - portability level - installability level - maintainability level - reliability level - reusability level - design - architecture - hours/days - time restriction - accuracy level.value - quality level - implementation - non functional requirements.other - requirements novelty level - requirements clarity level - processing requirements - requirements volatility level - project domain.other - project.infrastructure - estimation entity.other - entity count - performance - time efficiency level - process efficiency level - number of entities estimated.value - not used - availability - availability level - security - security level - reliability - robustness level - user case points - testing - testability level - maintenance'
$ws.Cells.Item(21, 1).Value = 'Usman - Britto_2016'
$ws.Cells.Item(21, 2).Value = 'distributed: far offshore - semi-distributed - distributed - distributed: distant onshore - co-located - relationship.location - centralized - distributed: near offshore - estimation entity.other - estimator & provider - estimator - distribution - estimation techniques.other - distributed: close onshore - estimate value(s)'
$ws.Cells.Item(22, 1).Value = 'Usman - Dasthi'
$ws.Cells.Item(22, 2).Value = 'analogy - analogy-based - expert judgement - expert judgment'
$ws.Cells.Item(23, 1).Value = 'Usman - Mendes'
$ws.Cells.Item(23, 2).Value = 'size.other - early size metric - late size metric - use case - functionality - not considered - validated theoretically - maintainability - task - implementation - non functional requirements.other - performance - considered - considered without any metric - solution-oriented metric - problem-oriented metric'

# --- Apply per-row font color/bold for column B based on target style groups ---
$ws.Cells.Item(2, 2).Font.Color = 12632139
$ws.Cells.Item(2, 2).Font.Bold = $false
$ws.Cells.Item(3, 2).Font.Color = 8676351
$ws.Cells.Item(3, 2).Font.Bold = $false
$ws.Cells.Item(4, 2).Font.Color = 8676351
$ws.Cells.Item(4, 2).Font.Bold = $true
$ws.Cells.Item(5, 2).Font.Color = 15442486
$ws.Cells.Item(5, 2).Font.Bold = $false
$ws.Cells.Item(6, 2).Font.Color = 15442486
$ws.Cells.Item(6, 2).Font.Bold = $true
$ws.Cells.Item(7, 2).Font.Color = 8676351
$ws.Cells.Item(7, 2).Font.Bold = $false
$ws.Cells.Item(8, 2).Font.Color = 12632139
$ws.Cells.Item(8, 2).Font.Bold = $false
$ws.Cells.Item(9, 2).Font.Color = 12632139
$ws.Cells.Item(9, 2).Font.Bold = $false
$ws.Cells.Item(10, 2).Font.Color = 8676351
$ws.Cells.Item(10, 2).Font.Bold = $false
$ws.Cells.Item(11, 2).Font.Color = 16737945
$ws.Cells.Item(11, 2).Font.Bold = $false
$ws.Cells.Item(12, 2).Font.Color = 15442486
$ws.Cells.Item(12, 2).Font.Bold = $false
$ws.Cells.Item(13, 2).Font.Color = 5689087
$ws.Cells.Item(13, 2).Font.Bold = $true
$ws.Cells.Item(14, 2).Font.Color = 5689087
$ws.Cells.Item(14, 2).Font.Bold = $false
$ws.Cells.Item(15, 2).Font.Color = 8676351
$ws.Cells.Item(15, 2).Font.Bold = $false
$ws.Cells.Item(16, 2).Font.Color = 16737945
$ws.Cells.Item(16, 2).Font.Bold = $false
$ws.Cells.Item(17, 2).Font.Color = 5689087
$ws.Cells.Item(17, 2).Font.Bold = $false
$ws.Cells.Item(18, 2).Font.Color = 16737945
$ws.Cells.Item(18, 2).Font.Bold = $false
$ws.Cells.Item(19, 2).Font.Color = 8676351
$ws.Cells.Item(19, 2).Font.Bold = $true
$ws.Cells.Item(20, 2).Font.Color = 8676351
$ws.Cells.Item(20, 2).Font.Bold = $false
$ws.Cells.Item(21, 2).Font.Color = 4235263
$ws.Cells.Item(21, 2).Font.Bold = $false
$ws.Cells.Item(22, 2).Font.Color = 16737945
$ws.Cells.Item(22, 2).Font.Bold = $false
$ws.Cells.Item(23, 2).Font.Color = 12632139
$ws.Cells.Item(23, 2).Font.Bold = $false
